$d = $word.ActiveDocument

# 1. Title heading + bold byline near the end (both occurrences identical text)
$d.Content.Find.Execute("Play Fire of Egypt Free - Merkur Slot Review", $true, $false, $false, $false, $false, $true, 1, $false, "Play Fire of Egypt for Free - Game Review", 2)

# 2. "What we like" bullet - graphics
$d.Content.Find.Execute("Effective graphics that capture the essence of ancient Egyptian culture", $true, $false, $false, $false, $false, $true, 1, $false, "Reflects the essence of ancient Egyptian culture and traditions", 2)

# 3. "What we like" bullet - winning combinations
$d.Content.Find.Execute("Numerous winning combinations for increased winnings", $true, $false, $false, $false, $false, $true, 1, $false, "Excellent quality symbols that are easily recognizable", 2)

# 4. "What we don't like" bullet - 4x4 grid
$d.Content.Find.Execute("Unfamiliar 4x4 grid may not appeal to all players", $true, $false, $false, $false, $false, $true, 1, $false, "Graphics are not innovative", 2)

# 5. "What we don't like" bullet - double or nothing
$d.Content.Find.Execute("High-risk double or nothing feature", $true, $false, $false, $false, $false, $true, 1, $false, "Double or nothing feature carries a high risk", 2)

# 6. Italic summary sentence near the end
$d.Content.Find.Execute("Read our review of Fire of Egypt, a Merkur slot game set in ancient Egypt. Play free and discover the mix of classic and modern gameplay mechanics.", $true, $false, $false, $false, $false, $true, 1, $false, "Read our review of Fire of Egypt and play this slot game for free. Explore ancient Egyptian culture and win big!", 2)
